$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("Settings")
$constants = $wb.Worksheets.Item("Constants")
$assets = $wb.Worksheets.Item("Assets")

# --- Settings sheet --------------------------------------------------
# Orchestrator queue name now points at the new HR onboarding process
$settings.Range("B2").Value = "HROnboarding"

# --- Assets sheet --------------------------------------------------
$assets.Range("A2").Value = "ExcelPath"
$assets.Range("B2").Value = "ExcelPath"

$assets.Range("A3").Value = "EmailCredentials"
$assets.Range("B3").Value = "EmailCredentials"

# --- Constants sheet ---------------------------------------------------
# New email notification settings
$constants.Range("A20").Value = "EmailSubject"

# Screenshot folder path now lives under the Data directory
$constants.Range("B5").Value = "Data\Exceptions_Screenshots"

$constants.Range("A21").Value = "EmailBody"
$constants.Range("B21").Value = "Hello, `nAn exception occurred during the automation process.  Please find the details below:`nException Source: @Source`nException Message: @Message`nA screenshot of the error has been attached for reference. Please see the attachment for more details.`nThank you and have a good day,`nRobot :)"
$constants.Range("B21").WrapText = $true
# Keep the default row height (mirrors the author's workbook, which did not autofit this row)
$constants.Rows.Item(21).RowHeight = 14.25

$constants.Range("B20").Value = "Automation Error!"

# --- Selection / active-sheet bookkeeping, mirroring the authored edit ---
$settings.Activate() | Out-Null
$settings.Range("B3").Select() | Out-Null

$assets.Activate() | Out-Null
$assets.Range("B8").Select() | Out-Null

$constants.Activate() | Out-Null
$constants.Range("B21").Select() | Out-Null
